# Change "python kinesis.py" to "python simple-producer.py" and move the
# "_GoBack" bookmark (which Word drops at the location of the last edit)
# from in front of the "Reference" heading to right after the newly
# edited text, inside the same list-item paragraph.

$d = $word.ActiveDocument

# 1. Update the command text.
$d.Content.Find.Execute("python kinesis.py", $true, $false, $false, $false, $false,
                         $true, 1, $false, "python simple-producer.py", 2)

# 2. Remove the bookmark from its old location (before the "Reference"
#    heading run) so it does not end up duplicated.
$d.Bookmarks("_GoBack").Delete()

# 3. Locate the end of the text we just inserted.
$found = $d.Content
$found.Find.Execute("python simple-producer.py")
$endPos = $found.End

# Placing a brand-new, *zero-length* bookmark exactly on a
# paragraph-end boundary (i.e. immediately before the paragraph mark)
# is mishandled by this host and silently relocates the bookmark to
# the very start of the document. Work around it: temporarily insert
# a throwaway character right after the text (so the insertion point
# is no longer sitting on that paragraph-end boundary), anchor the
# bookmark there, then delete the throwaway character again. The
# already-created bookmark correctly collapses back onto the boundary
# once the extra character is removed.
$insPoint = $d.Range($endPos, $endPos)
$insPoint.InsertAfter("X")

$bmRange = $d.Range($endPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$dummy = $d.Range($endPos, $endPos + 1)
$dummy.Delete()
